$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new row before the old row 116 (the "※4/8..." note row), which
# pushes that note row down to row 117 and gives the freshly inserted row
# 116 the same per-column number formats as the row above it (115).
$ws.Rows.Item(116).Insert() | Out-Null

# Populate the new data row 116 with the new day's figures.
$ws.Cells.Item(116, 1).Value = 43971
$ws.Cells.Item(116, 2).Value = 175
$ws.Cells.Item(116, 3).Value = 38346
$ws.Cells.Item(116, 4).Value = 46
$ws.Cells.Item(116, 5).Value = 7728

# Update the named Print_Area to extend through the new last row (117).
$n = $wb.Names.Item(1)
$n.RefersTo = "=相談件数!`$A`$1:`$E`$117"

# Move the cursor/selection to the new bottom-most cell, mirroring where the
# author's cursor ended up when they saved the workbook.
$ws.Range("B117").Select() | Out-Null

Write-Output "Updated print area and appended row for 2020-05-20; note row now at 117."
